# Convention change to support multi-axle vehicles:
# Add a new "Truck_Amandla" worksheet, cloned from "Bus_Makhulu", with
# updated driver-position offsets and its own identifying label.

$wb = $excel.ActiveWorkbook

# 1) Clone the Bus_Makhulu sheet (keeps its layout, styles, conditional
#    formatting, formulas, etc.) and place the copy at the end of the
#    workbook's tab strip.
$source = $wb.Worksheets.Item("Bus_Makhulu")
$source.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Truck_Amandla"

# 1b) Re-create the "class" conditional formatting rule on the clone so it
#     gets its own dxf style entry (same light-yellow fill as the other
#     sheets) instead of silently sharing the source sheet's entry.
$cfRange = $newSheet.Range("A4:B4")
$cfRange.FormatConditions.Delete()
$newRule = $cfRange.FormatConditions.Add(1, 3, '"class"')
$newRule.Interior.Color = 13431551

# 2) Update the sheet's own identifying label (H3) to match its new name.
$newSheet.Range("H3").Value = "Truck_Amandla"

# 3) Updated driver seat offsets for the multi-axle truck convention.
$newSheet.Range("F5").Value = -1.3230999999999999
$newSheet.Range("G5").Value = 0.55801299999999998
$newSheet.Range("H5").Value = 2.3923999999999999

# 4) Make the new sheet the active / selected tab, matching where the
#    author left their cursor.
$newSheet.Activate()
$newSheet.Range("N10").Select()
